# This script applies updated NATMI TPM-derived values to the
# "Ptn-Ptprs" ligand-receptor results sheet, reflecting a re-run of the
# analysis scripts with new TPM-normalized expression values
# (commit: "update scripts wuth new tpm").
#
# All changed cells hold literal numeric results produced by the
# upstream Python/NATMI pipeline (no formulas in this sheet), so the
# values are written directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07024999999999999
$ws.Range("H2").Value = 0.21075
$ws.Range("I2").Value = 0.005236595731231519
$ws.Range("J2").Value = 0.005236595731231519
$ws.Range("M2").Value = 1.660421
$ws.Range("N2").Value = 4.981262999999999
$ws.Range("O2").Value = 0.03714789785507311
$ws.Range("P2").Value = 0.03714789785507311
$ws.Range("Q2").Value = 0.11664457525
$ws.Range("R2").Value = 1.04980117725
$ws.Range("S2").Value = 0.0001945285233321004
$ws.Range("T2").Value = 0.0001945285233321004

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.07024999999999999
$ws.Range("H3").Value = 0.21075
$ws.Range("I3").Value = 0.005236595731231519
$ws.Range("J3").Value = 0.005236595731231519
$ws.Range("O3").Value = 0.5631392661118858
$ws.Range("P3").Value = 0.5631392661118859
$ws.Range("Q3").Value = 1.768259963416666
$ws.Range("R3").Value = 15.91433967075
$ws.Range("S3").Value = 0.002948932677010351
$ws.Range("T3").Value = 0.002948932677010352

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07024999999999999
$ws.Range("H4").Value = 0.21075
$ws.Range("I4").Value = 0.005236595731231519
$ws.Range("J4").Value = 0.005236595731231519
$ws.Range("M4").Value = 17.866195
$ws.Range("N4").Value = 53.598585
$ws.Range("O4").Value = 0.399712836033041
$ws.Range("P4").Value = 0.399712836033041
$ws.Range("Q4").Value = 1.25510019875
$ws.Range("R4").Value = 11.29590178875
$ws.Range("S4").Value = 0.002093134530889067
$ws.Range("T4").Value = 0.002093134530889067

# Row 5
$ws.Range("I5").Value = 0.4287876899474159
$ws.Range("J5").Value = 0.4287876899474159
$ws.Range("M5").Value = 1.660421
$ws.Range("N5").Value = 4.981262999999999
$ws.Range("O5").Value = 0.03714789785507311
$ws.Range("P5").Value = 0.03714789785507311
$ws.Range("Q5").Value = 9.551197100827665
$ws.Range("R5").Value = 85.960773907449
$ws.Range("S5").Value = 0.01592856130767937
$ws.Range("T5").Value = 0.01592856130767937

# Row 6
$ws.Range("I6").Value = 0.4287876899474159
$ws.Range("J6").Value = 0.4287876899474159
$ws.Range("O6").Value = 0.5631392661118858
$ws.Range("P6").Value = 0.5631392661118859
$ws.Range("S6").Value = 0.2414671850347986
$ws.Range("T6").Value = 0.2414671850347987

# Row 7
$ws.Range("I7").Value = 0.4287876899474159
$ws.Range("J7").Value = 0.4287876899474159
$ws.Range("M7").Value = 17.866195
$ws.Range("N7").Value = 53.598585
$ws.Range("O7").Value = 0.399712836033041
$ws.Range("P7").Value = 0.399712836033041
$ws.Range("Q7").Value = 102.7712549328283
$ws.Range("R7").Value = 924.9412943954551
$ws.Range("S7").Value = 0.1713919436049379
$ws.Range("T7").Value = 0.1713919436049379

# Row 8
$ws.Range("G8").Value = 7.592679666666666
$ws.Range("H8").Value = 22.778039
$ws.Range("I8").Value = 0.5659757143213526
$ws.Range("J8").Value = 0.5659757143213525
$ws.Range("M8").Value = 1.660421
$ws.Range("N8").Value = 4.981262999999999
$ws.Range("O8").Value = 0.03714789785507311
$ws.Range("P8").Value = 0.03714789785507311
$ws.Range("Q8").Value = 12.60704476480633
$ws.Range("R8").Value = 113.463402883257
$ws.Range("S8").Value = 0.02102480802406165
$ws.Range("T8").Value = 0.02102480802406164

# Row 9
$ws.Range("G9").Value = 7.592679666666666
$ws.Range("H9").Value = 22.778039
$ws.Range("I9").Value = 0.5659757143213526
$ws.Range("J9").Value = 0.5659757143213525
$ws.Range("O9").Value = 0.5631392661118858
$ws.Range("P9").Value = 0.5631392661118859
$ws.Range("Q9").Value = 191.1150387133732
$ws.Range("R9").Value = 1720.035348420359
$ws.Range("S9").Value = 0.3187231484000769
$ws.Range("T9").Value = 0.3187231484000769

# Row 10
$ws.Range("G10").Value = 7.592679666666666
$ws.Range("H10").Value = 22.778039
$ws.Range("I10").Value = 0.5659757143213526
$ws.Range("J10").Value = 0.5659757143213525
$ws.Range("M10").Value = 17.866195
$ws.Range("N10").Value = 53.598585
$ws.Range("O10").Value = 0.399712836033041
$ws.Range("P10").Value = 0.399712836033041
$ws.Range("Q10").Value = 135.6522954972017
$ws.Range("R10").Value = 1220.870659474815
$ws.Range("S10").Value = 0.2262277578972141
$ws.Range("T10").Value = 0.2262277578972141
